$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "culture_collection" attribute (column U) is being removed entirely
# again (re-deleted per INSDC2017 review, see commit message).
#
# Deleting the column shifts the header/data cells (and shared strings)
# one place to the left automatically, but cell comments stay anchored to
# their original column letters in this engine, so after the column
# delete we re-home every comment from V15..BW15 down to U15..BV15 by
# hand (overwriting each comment with the text that used to sit one
# column to its right), then remove the now-duplicated trailing comment.

$ws.Columns("U:U").Delete()

[void]($ws.Range("U15").Comment.Text('date of most recent douche'))
[void]($ws.Range("V15").Comment.Text('Traits like antibiotic resistance/xenobiotic degration phenotypes/converting phage genes'))
[void]($ws.Range("W15").Comment.Text('Estimated size of genome'))
[void]($ws.Range("X15").Comment.Text('ethnicity of the subject'))
[void]($ws.Range("Y15").Comment.Text('Plasmids that have significance phenotypic consequence'))
[void]($ws.Range("Z15").Comment.Text('history of gynecological disorders; can include multiple disorders'))
[void]($ws.Range("AA15").Comment.Text('Health or disease status of sample at time of collection'))
[void]($ws.Range("AB15").Comment.Text('Age of host at the time of sampling'))
[void]($ws.Range("AC15").Comment.Text('body mass index of the host, calculated as weight/(height)squared'))
[void]($ws.Range("AD15").Comment.Text('substance produced by the host, e.g. stool, mucus, where the sample was obtained from'))
[void]($ws.Range("AE15").Comment.Text('core body temperature of the host when sample was collected'))
[void]($ws.Range("AF15").Comment.Text('type of diet depending on the sample for animals omnivore, herbivore etc., for humans high-fat, meditteranean etc.; can include multiple diet types'))
[void]($ws.Range("AG15").Comment.Text('Name of relevant disease, e.g. Salmonella gastroenteritis. For the controlled vocabulary, please see Human Disease Ontology, http://bioportal.bioontology.org/ontologies/1009 or MeSH, http://www.ncbi.nlm.nih.gov/mesh'))
[void]($ws.Range("AH15").Comment.Delete())
[void]($ws.Range("AI15").Comment.Delete())
[void]($ws.Range("AJ15").Comment.Text('the height of subject'))
[void]($ws.Range("AK15").Comment.Text('content of last meal and time since feeding; can include multiple values'))
[void]($ws.Range("AL15").Comment.Text('most frequent job performed by subject'))
[void]($ws.Range("AM15").Comment.Delete())
[void]($ws.Range("AN15").Comment.Text('resting pulse of the host, measured as beats per minute'))
[void]($ws.Range("AO15").Comment.Text('Gender or physical sex of the host'))
[void]($ws.Range("AP15").Comment.Text('a unique identifier by which each subject can be referred to, de-identified, e.g. #131'))
[void]($ws.Range("AQ15").Comment.Text('NCBI taxonomy ID of the host, e.g. 9606'))
[void]($ws.Range("AR15").Comment.Text('Type of tissue the initial sample was taken from. Controlled vocabulary, http://bioportal.bioontology.org/ontologies/1005'))
[void]($ws.Range("AS15").Comment.Text('total mass of the host at collection, the unit depends on host'))
[void]($ws.Range("AT15").Comment.Text('whether subject had hormone replacement theraphy, and if yes start date'))
[void]($ws.Range("AU15").Comment.Text('specification of whether hysterectomy was performed'))
[void]($ws.Range("AV15").Comment.Text('can include multiple medication codes'))
[void]($ws.Range("AW15").Comment.Text('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.'))
[void]($ws.Range("AX15").Comment.Text('A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html'))
[void]($ws.Range("AY15").Comment.Text('whether full medical history was collected'))
[void]($ws.Range("AZ15").Comment.Text('date of most recent menstruation'))
[void]($ws.Range("BA15").Comment.Text('date of onset of menopause'))
[void]($ws.Range("BB15").Comment.Text('any other measurement performed or parameter collected, that is not listed here'))
[void]($ws.Range("BC15").Comment.Text('total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts'))
[void]($ws.Range("BD15").Comment.Text('oxygenation status of sample'))
[void]($ws.Range("BE15").Comment.Text('To what is the entity pathogenic'))
[void]($ws.Range("BF15").Comment.Text('type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types'))
[void]($ws.Range("BG15").Comment.Text('date due of pregnancy'))
[void]($ws.Range("BH15").Comment.Text('Aerobic or anaerobic'))
[void]($ws.Range("BI15").Comment.Text('Method or device employed for collecting sample'))
[void]($ws.Range("BJ15").Comment.Text('Processing applied to the sample during or after isolation'))
[void]($ws.Range("BK15").Comment.Text('salinity of sample, i.e. measure of total salt concentration'))
[void]($ws.Range("BL15").Comment.Text('Amount or size of sample (volume, mass or area) that was collected'))
[void]($ws.Range("BM15").Comment.Text('duration for which sample was stored'))
[void]($ws.Range("BN15").Comment.Text('location at which sample was stored, usually name of a specific freezer/room'))
[void]($ws.Range("BO15").Comment.Text('temperature at which sample was stored, e.g. -80'))
[void]($ws.Range("BP15").Comment.Text('volume (mL) or weight (g) of sample processed for DNA extraction'))
[void]($ws.Range("BQ15").Comment.Text('current sexual partner and frequency of sex'))
[void]($ws.Range("BR15").Comment.Text('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.'))
[void]($ws.Range("BS15").Comment.Text('Information about the genetic distinctness of the lineage (eg., biovar, serovar)'))
[void]($ws.Range("BT15").Comment.Text('temperature of the sample at time of sampling'))
[void]($ws.Range("BU15").Comment.Text('Feeding position in food chain (eg., chemolithotroph)'))
[void]($ws.Range("BV15").Comment.Text('history of urogenital disorders, can include multiple disorders'))
[void]($ws.Range("BW15").Comment.Delete())
